$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the obsolete rows (old rows 7-13) ---
$ws.Range("A7:A13").EntireRow.Delete()

# --- Drop all existing hyperlinks; they'll be re-created below for the
#     rows that survive, pointing at the refreshed URLs ---
$ws.Range("A1").Hyperlinks.Delete()

# --- Row 2 ---
$ws.Range("A2").Value = "2025-11-06 06:27:35"
$ws.Range("B2").Value = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5427956"
$ws.Range("G2").Value = 310
$ws.Range("H2").Value = "🔥AI,Ai"

# --- Row 3 ---
$ws.Range("A3").Value = "2025-11-06 06:27:35"
$ws.Range("B3").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G3").Value = 243
$ws.Range("H3").Value = "🔥API ◆ツール"

# --- Row 4 ---
$ws.Range("A4").Value = "2025-11-06 06:27:35"
$ws.Range("B4").Value = "【急募】GitHub管理のBootstrapサイト移行作業依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5428337"
$ws.Range("G4").Value = 58
$ws.Range("H4").Value = "◇サイト"

# --- Row 5 ---
$ws.Range("A5").Value = "2025-11-06 06:27:35"
$ws.Range("B5").Value = "勤怠管理システムの改修依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5428278"
$ws.Range("G5").Value = 53
$ws.Range("H5").Value = "◇管理"

# --- Row 6 ---
$ws.Range("A6").Value = "2025-11-06 06:27:35"
$ws.Range("B6").Value = "【急募】WEB会計アプリ機能修正!納期11/09希望"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5428124"
$ws.Range("G6").Value = 38
$ws.Range("H6").Value = "◇アプリ"

# --- Re-create hyperlinks for the URL column on the surviving rows ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5427956")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5428337")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5428278")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5428124")

# Hyperlinks.Add() swaps in a freshly-minted "Hyperlink" style xf; put the
# cells back on the sheet's original Hyperlink style so formatting matches.
$ws.Range("F2:F6").Style = "Hyperlink"

# --- Column width tweaks ---
# ColumnWidth (characters) round-trips into the stored OOXML <col width>
# with a +5/6 offset in this engine, so back it out to land on the exact
# target widths (51 and 12).
$ws.Columns.Item(2).ColumnWidth = 50.166666666666664
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666
